$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "[[`"Steel cylinder`", 2.4821260976820004e-05], [`"nitrous oxide production`", 0.010968005818038482], [`"Use phase basecase`", 0.0], [`"EoL 50L cylinder - FU`", 0.0015934554481004321]]"
$ws.Range("B2").Value = "[[`"Steel cylinder`", 0.005456229514756893], [`"nitrous oxide production`", 2.7735767368862367], [`"Use phase basecase`", 227.34409189224243], [`"EoL 50L cylinder - FU`", 1.3797549866934615]]"
$ws.Range("C2").Value = "[[`"Steel cylinder`", 0.02689681145503979], [`"nitrous oxide production`", 15.635308294302359], [`"Use phase basecase`", 0.0], [`"EoL 50L cylinder - FU`", 2.5968806464953906]]"
$ws.Range("D2").Value = "[[`"Steel cylinder`", 0.06471586600422786], [`"nitrous oxide production`", 46.31028189903245], [`"Use phase basecase`", 0.0], [`"EoL 50L cylinder - FU`", 6.953517004338244]]"
$ws.Range("E2").Value = "[[`"Steel cylinder`", 4.164145755153704e-06], [`"nitrous oxide production`", 0.0009169092173445685], [`"Use phase basecase`", 0.0], [`"EoL 50L cylinder - FU`", 0.0001231153606355252]]"
$ws.Range("F2").Value = "[[`"Steel cylinder`", 5.439366377934605e-06], [`"nitrous oxide production`", 0.005230909219588117], [`"Use phase basecase`", 0.0], [`"EoL 50L cylinder - FU`", 0.0006537394847833525]]"
$ws.Range("G2").Value = "[[`"Steel cylinder`", 5.1440088432858044e-05], [`"nitrous oxide production`", 0.029899230463375122], [`"Use phase basecase`", 0.0], [`"EoL 50L cylinder - FU`", 0.004205256078726926]]"
$ws.Range("H2").Value = "[[`"Steel cylinder`", 8.120702205798173e-11], [`"nitrous oxide production`", 1.3815806307562065e-09], [`"Use phase basecase`", 0.0], [`"EoL 50L cylinder - FU`", 2.4542745026166035e-10]]"
$ws.Range("I2").Value = "[[`"Steel cylinder`", 6.725000501277105e-11], [`"nitrous oxide production`", 2.7210021087733674e-08], [`"Use phase basecase`", 1.3361670158765134e-08], [`"EoL 50L cylinder - FU`", 5.05522698476935e-09]]"
$ws.Range("J2").Value = "[[`"Steel cylinder`", 3.451402200181894e-05], [`"nitrous oxide production`", 0.42490983327819815], [`"Use phase basecase`", 0.0], [`"EoL 50L cylinder - FU`", 0.052218194515324355]]"
$ws.Range("K2").Value = "[[`"Steel cylinder`", 0.05341058184255911], [`"nitrous oxide production`", 28.792690550160245], [`"Use phase basecase`", 0.0], [`"EoL 50L cylinder - FU`", 4.553405133415671]]"
$ws.Range("L2").Value = "[[`"Steel cylinder`", -6.260798244751308e-08], [`"nitrous oxide production`", 3.528425210128439e-05], [`"Use phase basecase`", 0.0], [`"EoL 50L cylinder - FU`", 5.124795525628333e-06]]"
$ws.Range("M2").Value = "[[`"Steel cylinder`", 1.232676817653973e-10], [`"nitrous oxide production`", 8.530261057304866e-08], [`"Use phase basecase`", 0.0], [`"EoL 50L cylinder - FU`", 1.043284139163744e-07]]"
$ws.Range("N2").Value = "[[`"Steel cylinder`", 5.476745429210785e-10], [`"nitrous oxide production`", 5.905655038504484e-08], [`"Use phase basecase`", 0.0], [`"EoL 50L cylinder - FU`", 1.604819459417503e-08]]"
$ws.Range("O2").Value = "[[`"Steel cylinder`", 2.334134230290092e-05], [`"nitrous oxide production`", 0.007059262000778557], [`"Use phase basecase`", 0.0], [`"EoL 50L cylinder - FU`", 0.0012673930560681216]]"
$ws.Range("P2").Value = "[[`"Steel cylinder`", 0.0015157174317935512], [`"nitrous oxide production`", 4.430496838846222], [`"Use phase basecase`", 0.0], [`"EoL 50L cylinder - FU`", 0.4699444626920632]]"
$ws.Range("A3").Value = "[[`"CDU updated`", 0.033699763481384294], [`"Steel cylinder`", 2.4821260976820004e-05], [`"nitrous oxide production`", 0.010968005818038482], [`"Use phase CDU`", 0.0011995761189071566], [`"EoL CDU`", 0.0012242722397637902], [`"EoL 50L cylinder - FU`", 0.0015934554481004321]]"
$ws.Range("B3").Value = "[[`"CDU updated`", 0.12819826622266475], [`"Steel cylinder`", 0.005456229514756893], [`"nitrous oxide production`", 2.7735767368862367], [`"Use phase CDU`", 43.445356377489226], [`"EoL CDU`", -0.00030258188038357025], [`"EoL 50L cylinder - FU`", 1.3797549866934615]]"
$ws.Range("C3").Value = "[[`"CDU updated`", 1.4755335896694415], [`"Steel cylinder`", 0.02689681145503979], [`"nitrous oxide production`", 15.635308294302359], [`"Use phase CDU`", 0.7412027382264046], [`"EoL CDU`", -0.004626877616344883], [`"EoL 50L cylinder - FU`", 2.5968806464953906]]"
$ws.Range("D3").Value = "[[`"CDU updated`", 2.094289083153275], [`"Steel cylinder`", 0.06471586600422786], [`"nitrous oxide production`", 46.31028189903245], [`"Use phase CDU`", 0.6692525250789959], [`"EoL CDU`", 0.0023623570522492655], [`"EoL 50L cylinder - FU`", 6.953517004338244]]"
$ws.Range("E3").Value = "[[`"CDU updated`", 0.00010946578693302614], [`"Steel cylinder`", 4.164145755153704e-06], [`"nitrous oxide production`", 0.0009169092173445685], [`"Use phase CDU`", 2.849641883337073e-05], [`"EoL CDU`", -6.603880166606268e-07], [`"EoL 50L cylinder - FU`", 0.0001231153606355252]]"
$ws.Range("F3").Value = "[[`"CDU updated`", 0.00012323944886451918], [`"Steel cylinder`", 5.439366377934605e-06], [`"nitrous oxide production`", 0.005230909219588117], [`"Use phase CDU`", 0.00027925130312526306], [`"EoL CDU`", -6.2441110170754856e-06], [`"EoL 50L cylinder - FU`", 0.0006537394847833525]]"
$ws.Range("G3").Value = "[[`"CDU updated`", 0.001291805268237013], [`"Steel cylinder`", 5.1440088432858044e-05], [`"nitrous oxide production`", 0.029899230463375122], [`"Use phase CDU`", 0.004351892777136212], [`"EoL CDU`", -9.307467647554836e-05], [`"EoL 50L cylinder - FU`", 0.004205256078726926]]"
$ws.Range("H3").Value = "[[`"CDU updated`", 3.23120762737068e-10], [`"Steel cylinder`", 8.120702205798173e-11], [`"nitrous oxide production`", 1.3815806307562065e-09], [`"Use phase CDU`", 3.112433080587344e-10], [`"EoL CDU`", -7.707922940565687e-12], [`"EoL 50L cylinder - FU`", 2.4542745026166035e-10]]"
$ws.Range("I3").Value = "[[`"CDU updated`", 4.527936227184927e-09], [`"Steel cylinder`", 6.725000501277105e-11], [`"nitrous oxide production`", 2.7210021087733674e-08], [`"Use phase CDU`", 1.1108940558000346e-08], [`"EoL CDU`", -2.91480127867933e-10], [`"EoL 50L cylinder - FU`", 5.05522698476935e-09]]"
$ws.Range("J3").Value = "[[`"CDU updated`", 0.022953633802310615], [`"Steel cylinder`", 3.451402200181894e-05], [`"nitrous oxide production`", 0.42490983327819815], [`"Use phase CDU`", 0.00034780368107774893], [`"EoL CDU`", 0.00011107874461068391], [`"EoL 50L cylinder - FU`", 0.052218194515324355]]"
$ws.Range("K3").Value = "[[`"CDU updated`", 1.4453593855555278], [`"Steel cylinder`", 0.05341058184255911], [`"nitrous oxide production`", 28.792690550160245], [`"Use phase CDU`", 30.062931903317136], [`"EoL CDU`", -0.02832465624823723], [`"EoL 50L cylinder - FU`", 4.553405133415671]]"
$ws.Range("L3").Value = "[[`"CDU updated`", 2.3204757355396628e-05], [`"Steel cylinder`", -6.260798244751308e-08], [`"nitrous oxide production`", 3.528425210128439e-05], [`"Use phase CDU`", 9.73441196893503e-06], [`"EoL CDU`", 7.408005945760737e-07], [`"EoL 50L cylinder - FU`", 5.124795525628333e-06]]"
$ws.Range("M3").Value = "[[`"CDU updated`", 6.412188285949878e-09], [`"Steel cylinder`", 1.232676817653973e-10], [`"nitrous oxide production`", 8.530261057304866e-08], [`"Use phase CDU`", 1.9450074604811147e-09], [`"EoL CDU`", 1.3108767544701095e-11], [`"EoL 50L cylinder - FU`", 1.043284139163744e-07]]"
$ws.Range("N3").Value = "[[`"CDU updated`", 2.9190976038666193e-08], [`"Steel cylinder`", 5.476745429210785e-10], [`"nitrous oxide production`", 5.905655038504484e-08], [`"Use phase CDU`", 1.3992546508923827e-08], [`"EoL CDU`", 7.644866224450964e-10], [`"EoL 50L cylinder - FU`", 1.604819459417503e-08]]"
$ws.Range("O3").Value = "[[`"CDU updated`", 0.002678400846522548], [`"Steel cylinder`", 2.334134230290092e-05], [`"nitrous oxide production`", 0.007059262000778557], [`"Use phase CDU`", 0.0008604984653242085], [`"EoL CDU`", 5.4636027418626265e-05], [`"EoL 50L cylinder - FU`", 0.0012673930560681216]]"
$ws.Range("P3").Value = "[[`"CDU updated`", 0.04610221142028905], [`"Steel cylinder`", 0.0015157174317935512], [`"nitrous oxide production`", 4.430496838846222], [`"Use phase CDU`", 0.026577952306826722], [`"EoL CDU`", 0.000897610415032529], [`"EoL 50L cylinder - FU`", 0.4699444626920632]]"
$ws.Range("A4").Value = "[[`"MDU updated`", 0.06909101812686011], [`"Steel cylinder`", 2.4821260976820004e-05], [`"nitrous oxide production`", 0.010968005818038482], [`"Use phase MDU`", 0.00045196796330177], [`"EoL MDU`", 0.007055424060427032], [`"EoL 50L cylinder - FU`", 0.0015934554481004321]]"
$ws.Range("B4").Value = "[[`"MDU updated`", 0.3591249304907266], [`"Steel cylinder`", 0.005456229514756893], [`"nitrous oxide production`", 2.7735767368862367], [`"Use phase MDU`", 43.404329830371644], [`"EoL MDU`", -0.0070986448114089525], [`"EoL 50L cylinder - FU`", 1.3797549866934615]]"
$ws.Range("C4").Value = "[[`"MDU updated`", 5.562379935751366], [`"Steel cylinder`", 0.02689681145503979], [`"nitrous oxide production`", 15.635308294302359], [`"Use phase MDU`", 0.2792652226980472], [`"EoL MDU`", -0.10380338196557717], [`"EoL 50L cylinder - FU`", 2.5968806464953906]]"
$ws.Range("D4").Value = "[[`"MDU updated`", 5.548337101274091], [`"Steel cylinder`", 0.06471586600422786], [`"nitrous oxide production`", 46.31028189903245], [`"Use phase MDU`", 0.2521563208261186], [`"EoL MDU`", -0.03427735018172776], [`"EoL 50L cylinder - FU`", 6.953517004338244]]"
$ws.Range("E4").Value = "[[`"MDU updated`", 0.0003663742164005576], [`"Steel cylinder`", 4.164145755153704e-06], [`"nitrous oxide production`", 0.0009169092173445685], [`"Use phase MDU`", 1.073668288198834e-05], [`"EoL MDU`", -9.489777176417246e-06], [`"EoL 50L cylinder - FU`", 0.0001231153606355252]]"
$ws.Range("F4").Value = "[[`"MDU updated`", 0.000361856412791819], [`"Steel cylinder`", 5.439366377934605e-06], [`"nitrous oxide production`", 0.005230909219588117], [`"Use phase MDU`", 0.00010521436758666906], [`"EoL MDU`", -4.8532739762006666e-05], [`"EoL 50L cylinder - FU`", 0.0006537394847833525]]"
$ws.Range("G4").Value = "[[`"MDU updated`", 0.0037610866718165645], [`"Steel cylinder`", 5.1440088432858044e-05], [`"nitrous oxide production`", 0.029899230463375122], [`"Use phase MDU`", 0.0016396759521870126], [`"EoL MDU`", -0.0007152441329165199], [`"EoL 50L cylinder - FU`", 0.004205256078726926]]"
$ws.Range("H4").Value = "[[`"MDU updated`", 1.0716849152239044e-09], [`"Steel cylinder`", 8.120702205798173e-11], [`"nitrous oxide production`", 1.3815806307562065e-09], [`"Use phase MDU`", 1.1726809313497473e-10], [`"EoL MDU`", -8.165715833980981e-11], [`"EoL 50L cylinder - FU`", 2.4542745026166035e-10]]"
$ws.Range("I4").Value = "[[`"MDU updated`", 1.1984131035265172e-08], [`"Steel cylinder`", 6.725000501277105e-11], [`"nitrous oxide production`", 2.7210021087733674e-08], [`"Use phase MDU`", 5.7744917398854014e-09], [`"EoL MDU`", -4.293254963228643e-09], [`"EoL 50L cylinder - FU`", 5.05522698476935e-09]]"
$ws.Range("J4").Value = "[[`"MDU updated`", 0.05186765582765432], [`"Steel cylinder`", 3.451402200181894e-05], [`"nitrous oxide production`", 0.42490983327819815], [`"Use phase MDU`", 0.00013104305669854347], [`"EoL MDU`", 0.0005209868928899672], [`"EoL 50L cylinder - FU`", 0.052218194515324355]]"
$ws.Range("K4").Value = "[[`"MDU updated`", 3.6900959224740673], [`"Steel cylinder`", 0.05341058184255911], [`"nitrous oxide production`", 28.792690550160245], [`"Use phase MDU`", 11.326902802631512], [`"EoL MDU`", -0.29824717604918993], [`"EoL 50L cylinder - FU`", 4.553405133415671]]"
$ws.Range("L4").Value = "[[`"MDU updated`", 7.762141348560289e-05], [`"Steel cylinder`", -6.260798244751308e-08], [`"nitrous oxide production`", 3.528425210128439e-05], [`"Use phase MDU`", 3.6676641708639582e-06], [`"EoL MDU`", -8.07769332375321e-07], [`"EoL 50L cylinder - FU`", 5.124795525628333e-06]]"
$ws.Range("M4").Value = "[[`"MDU updated`", 1.8962760185311458e-08], [`"Steel cylinder`", 1.232676817653973e-10], [`"nitrous oxide production`", 8.530261057304866e-08], [`"Use phase MDU`", 7.328264098056545e-10], [`"EoL MDU`", 2.667194990667454e-11], [`"EoL 50L cylinder - FU`", 1.043284139163744e-07]]"
$ws.Range("N4").Value = "[[`"MDU updated`", 6.342169280380604e-08], [`"Steel cylinder`", 5.476745429210785e-10], [`"nitrous oxide production`", 5.905655038504484e-08], [`"Use phase MDU`", 5.272014545275251e-09], [`"EoL MDU`", 4.049234970385121e-09], [`"EoL 50L cylinder - FU`", 1.604819459417503e-08]]"
$ws.Range("O4").Value = "[[`"MDU updated`", 0.005810953630291378], [`"Steel cylinder`", 2.334134230290092e-05], [`"nitrous oxide production`", 0.007059262000778557], [`"Use phase MDU`", 0.000324212638670387], [`"EoL MDU`", 0.0002787249201572001], [`"EoL 50L cylinder - FU`", 0.0012673930560681216]]"
$ws.Range("P4").Value = "[[`"MDU updated`", 0.1236924037638626], [`"Steel cylinder`", 0.0015157174317935512], [`"nitrous oxide production`", 4.430496838846222], [`"Use phase MDU`", 0.01001385638102837], [`"EoL MDU`", 0.004563970418034172], [`"EoL 50L cylinder - FU`", 0.4699444626920632]]"

Write-Output "Updated 48 cells"
